# "added more git commands"
# Append two new git-command rows (command + description) to the bottom
# of the "Git" command-reference sheet, right after the existing list
# that currently ends at row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Git")
$ws.Activate()

$ws.Range("A13").Value = "git commit -a"
$ws.Range("C13").Value = "commits all modified files to local repository"

$ws.Range("A14").Value = "git pull --recurse-submodules"
$ws.Range("C14").Value = "pulls everything from the remote repository and submodules and places it in the local repository"

# Leave the selection where the author would naturally end up after
# typing the last entry (one row below the newly-added data).
$ws.Range("C15").Select()
